$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 57, shifting existing rows 57-74 down to 58-75.
$ws.Rows.Item(57).Insert()

# Populate the newly inserted row 57 with the new weekly price record.
$ws.Cells.Item(57, 1).Value = 2
$ws.Cells.Item(57, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(57, 3).Value = "Coquimbo"
$ws.Cells.Item(57, 4).Value = 44839
$ws.Cells.Item(57, 5).Value = 4
$ws.Cells.Item(57, 6).Value = 100112026
$ws.Cells.Item(57, 7).Value = "Haba"
$ws.Cells.Item(57, 8).Value = "Sin especificar"
$ws.Cells.Item(57, 9).Value = "Primera"
$ws.Cells.Item(57, 10).Value = 1200
$ws.Cells.Item(57, 11).Value = 5000
$ws.Cells.Item(57, 12).Value = 6000
$ws.Cells.Item(57, 13).Value = 5500
$ws.Cells.Item(57, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(57, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(57, 16).Value = 220
$ws.Cells.Item(57, 17).Value = 25
$ws.Cells.Item(57, 18).Value = "Hortaliza"
